$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 281, shifting existing rows 281:333 down to 282:334.
$ws.Rows("281:281").Insert()

# Populate the newly inserted row 281 with the new weekly data point.
$ws.Range("A281").Value = 3
$ws.Range("B281").Value = "Femacal de La Calera"
$ws.Range("C281").Value = "Coquimbo"
$ws.Range("D281").Value = 45244
$ws.Range("E281").Value = 5
$ws.Range("F281").Value = 100112026
$ws.Range("G281").Value = "Haba"
$ws.Range("H281").Value = "Sin especificar"
$ws.Range("I281").Value = "Primera"
$ws.Range("J281").Value = 50
$ws.Range("K281").Value = 9000
$ws.Range("L281").Value = 9000
$ws.Range("M281").Value = 9000
$ws.Range("N281").Value = "$/saco 25 kilos"
$ws.Range("O281").Value = "Provincia de Quillota"
$ws.Range("P281").Value = 360
$ws.Range("Q281").Value = 25
$ws.Range("R281").Value = "Hortaliza"
